$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.519801735877991
$ws.Range("B1").Value = 4.256080150604248
$ws.Range("C1").Value = 3.492135524749756
$ws.Range("D1").Value = 1.448733568191528
$ws.Range("E1").Value = 0.964720606803894
